$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 1074
$ws.Range("F6").Value = 2768
$ws.Range("F8").Value = 1348
$ws.Range("F11").Value = 968
$ws.Range("F12").Value = 1233
$ws.Range("F14").Value = 123
$ws.Range("F15").Value = 772
$ws.Range("F16").Value = 813
$ws.Range("F18").Value = 569
$ws.Range("F19").Value = 1155
$ws.Range("F21").Value = 688
$ws.Range("F22").Value = 629
$ws.Range("F23").Value = 240
$ws.Range("F24").Value = 335
$ws.Range("F26").Value = 706
$ws.Range("F27").Value = 705
$ws.Range("F28").Value = 8059
$ws.Range("F34").Value = 207
$ws.Range("F35").Value = 1673
$ws.Range("F37").Value = 169
$ws.Range("F39").Value = 155
$ws.Range("F42").Value = 167

$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 73
$ws.Range("F13").Value = 53
$ws.Range("F17").Value = 231

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 771

$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1074
$ws.Range("F5").Value = 2768
$ws.Range("F7").Value = 1348
$ws.Range("F9").Value = 1233
$ws.Range("F10").Value = 307
$ws.Range("F12").Value = 123
$ws.Range("F13").Value = 772
$ws.Range("F16").Value = 813
$ws.Range("F18").Value = 569
$ws.Range("F19").Value = 1155
$ws.Range("F21").Value = 73
$ws.Range("F22").Value = 688
$ws.Range("F23").Value = 629
$ws.Range("F24").Value = 240
$ws.Range("F25").Value = 335
$ws.Range("F27").Value = 705
$ws.Range("F28").Value = 8059
$ws.Range("F32").Value = 207
$ws.Range("F33").Value = 1673
$ws.Range("F36").Value = 155
$ws.Range("F37").Value = 53
$ws.Range("F38").Value = 53
$ws.Range("F42").Value = 167
